# DevOps_Tracker.xlsx - "adding more files to repo"
#
# On the "MonitoringTools" worksheet:
#  - The "Demo on ElasticSearch, Logstash and Kibana (ELK)" entry (row 1)
#    gains an extra reference link (webkid.io) and grows taller.
#  - A new row is appended at the bottom that repeats the same "Demo on
#    ElasticSearch, Logstash and Kibana (ELK):" header as a (currently
#    empty) placeholder entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MonitoringTools")

# Row 1: append the new link to the existing ELK demo entry.
$ws.Range("A1").Value = "Demo on ElasticSearch, Logstash and Kibana (ELK):`nhttps://www.youtube.com/watch?v=Kqs7UcCJquM`nhttp://blog.webkid.io/visualize-datasets-with-elk/`n"
$ws.Rows.Item(1).RowHeight = 57.6

# New row 5: placeholder entry with just the ELK demo header text.
$ws.Range("A5").Value = "Demo on ElasticSearch, Logstash and Kibana (ELK):`n"
$ws.Rows.Item(5).RowHeight = 43.2

# Leave the freshly-edited cell selected.
$ws.Range("A1").Select()
